# "data build and updatae" - the 2021 row (the first data row under the
# header) was removed from the trade-data table, so every subsequent
# year's row shifts up by one. Reproduce this the way a person would in
# Excel: select the entire row 2 (A2:XFD2) and delete it, shifting the
# remaining rows up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = $ws.Range("A2:XFD2")
[void]$row2.Select()
$row2.EntireRow.Delete()
